# "get data update model" - append the latest monthly fuel-price data rows
# to Sheet1 (Jun-2017 .. Oct-2017), then leave the workbook positioned the
# way the author left it: Sheet1 scrolled down to the new rows with B259(+)
# selected, and Sheet2 the active/selected tab.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

# Make sure we are adding onto Sheet1.
$ws1.Activate() | Out-Null

# Shared-string columns used by every data row.
$sourceUrl = "http://www.energy.gov.za/files/petroleum_frame.html"
$updateNo = "no"

# New rows: Date(serial), BFP, Crude_Oil, Fuel_Tax, Customs_Excise,
# Transport_Cost, Wholesale_Margin, Retail_Margin, Slate_Levy, Delivery_Cost,
# DSML, Equalization_Fund_Levy, Petroleum_Products_Levy, Full_95_Coast,
# Basic_Fuel_Price, Secondary_Storage, Secondary_Distribution,
# Oil_price_rands, U95_BFP, U95_OPR
$rows = @(
    @(42887, 1305, 572.66999999999996, 13.258699999999999, 50.06,  572.66999999999996, 315, 4, 0, 163, 41.5, 0.33, 35.6, 17.899999999999999, 17.3, 176.4, 0, 0, 10),
    @(42917, 1237, 504.67,             12.8756,             46.48, 504.67,             315, 4, 0, 163, 41.5, 0.33, 35.6, 17.899999999999999, 17.3, 176.4, 0, 0, 10),
    @(42948, 1256, 523.66999999999996, 13.15,               48.23, 523.66999999999996, 315, 4, 0, 163, 41.5, 0.33, 35.6, 17.899999999999999, 17.3, 176.4, 0, 0, 10),
    @(42979, 1323, 586.07000000000005, 13.215199999999999, 51.66, 586.07000000000005, 315, 4, 0, 163, 41.5, 0.33, 35.6, 17.899999999999999, 17.3, 181,   0, 0, 10),
    @(43009, 1352, 615.07000000000005, 13.1317,             55.99, 615.07000000000005, 315, 4, 0, 163, 41.5, 0.33, 35.6, 17.899999999999999, 17.3, 181,   0, 0, 10)
)

$startRow = 260
for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $startRow + $i
    $vals = $rows[$i]
    for ($c = 1; $c -le $vals.Count; $c++) {
        $ws1.Cells.Item($r, $c).Value = $vals[$c - 1]
    }
    $ws1.Cells.Item($r, 20).Value = $sourceUrl
    $ws1.Cells.Item($r, 21).Value = $updateNo
}

# Leave the cursor where the author left it: B259, the row just above the
# newly appended data.
$ws1.Range("B259").Select() | Out-Null

# The author ended the session back on Sheet2 (it is the selected/active tab
# on save).
$ws2.Activate() | Out-Null
